# "permutation problems with pick and don't pick way"
# Adds 3 new problem rows (11-13) to the "General problems" sheet for the
# new Permutation-related problems, and updates the active sheet/selection.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("General problems")

# Reuse the existing date-formatted style (same as A2:A10) for the new date
# cells, and the existing wrap-text style (same as G10) for the long B13
# description, by copying formats instead of assigning ad-hoc NumberFormat/
# WrapText (which would create brand-new style entries in styles.xml).
$ws1.Range("A10").Copy()
$ws1.Range("A11:A13").PasteSpecial(-4122)

$ws1.Range("G10").Copy()
$ws1.Range("B13").PasteSpecial(-4122)

# --- Row 11: Permutation of adding spaces in between of a string ---
$ws1.Range("A11").Value = 44399
$ws1.Range("B11").Value = "Permutation of adding spaces in between of a string"
$ws1.Range("C11").Value = "Done - D"
$ws1.Range("D11").Value = "Easy"
$ws1.Range("E11").Value = "PermutationWithSpaces"
$ws1.Range("F11").Value = "Recursion"

# --- Row 12: Permutation of changing case in a string ---
$ws1.Range("A12").Value = 44399
$ws1.Range("B12").Value = "Permutation of changing case in a string"
$ws1.Range("C12").Value = "Done - D"
$ws1.Range("D12").Value = "Easy"
$ws1.Range("E12").Value = "PermutationWithSpaces"
$ws1.Range("F12").Value = "Recursion"

# --- Row 13: Permutation of a string with both letters and digits ---
$ws1.Range("A13").Value = 44399
$ws1.Range("C13").Value = "Done - D"
$ws1.Range("D13").Value = "Easy"
$ws1.Range("E13").Value = "LetterCasePermutation"
$ws1.Range("B13").Value = "Permutation of string containing both alphabets and numbers. Case of alphabets will be toggled and numbers will be copied as it is."
$ws1.Range("F13").Value = "Recursion"
$ws1.Rows.Item(13).RowHeight = 43.5

# --- Update selections on the other two touched sheets ---
$ws7 = $wb.Worksheets.Item("Sliding Window")
$ws7.Range("A7").Select()

$ws8 = $wb.Worksheets.Item("Sorting")
$ws8.Range("A1").Select()

# --- "General problems" ends up the active/visible tab, selection on B13 ---
$ws1.Activate()
$ws1.Range("B13").Select()
